$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")
$wsAssets    = $wb.Worksheets.Item("Assets")

# ---------------------------------------------------------------------------
# Assets sheet - SCM_URL asset
# ---------------------------------------------------------------------------
$wsAssets.Activate()
$wsAssets.Range("A2").Value = "SCM_URL"
$wsAssets.Range("B2").Value = "SCM_URL"

# ---------------------------------------------------------------------------
# Settings sheet - the orchestrator queue name was renamed for the
# Purchase Order flow
# ---------------------------------------------------------------------------
$wsSettings.Activate()
$wsSettings.Range("B2").Value = "QUE_PONumbers"
$wsSettings.Range("B2").Style = "Normal"

# ---------------------------------------------------------------------------
# Assets sheet - remaining PO login assets
# ---------------------------------------------------------------------------
$wsAssets.Activate()
$wsAssets.Range("A4").Value = "POLogin_Username"
$wsAssets.Range("B4").Value = "POLogin_Username"

$wsAssets.Range("A5").Value = "POLogin_Password"
$wsAssets.Range("B5").Value = "POLogin_Password"

$wsAssets.Range("A3").Value = "POLogin_URL"
$wsAssets.Range("B3").Value = "POLogin_URL"

# ---------------------------------------------------------------------------
# Settings sheet - new email notification + state assignment lookup rows
# ---------------------------------------------------------------------------
$wsSettings.Activate()

$wsSettings.Range("A6").Value = "Email Recipient"
$wsSettings.Range("A7").Value = "Email Subject"
$wsSettings.Range("A8").Value = "Email Body"
$wsSettings.Range("A9").Value = "Email Attachment"

$wsSettings.Range("B6").Value = "purni.work@gmail.com"
$wsSettings.Range("B7").Value = "PO Submission - Screen Shot"
$wsSettings.Range("B8").Value = "Please view the attached screen Shot"
$wsSettings.Range("B9").Value = "Data\Output\ScreenShot.png"

$wsSettings.Range("A10").Value = "StateAssignments"
$wsSettings.Range("B10").Value = "Data\Input\StateAssignments.xlsx"

# The recipient e-mail address is turned into a mailto hyperlink
$wsSettings.Hyperlinks.Add($wsSettings.Range("B6"), "mailto:purni.work@gmail.com", "", "", "purni.work@gmail.com")

$wsSettings.Range("A9").Select()

# ---------------------------------------------------------------------------
# Assets sheet - final PO lookup URL asset
# ---------------------------------------------------------------------------
$wsAssets.Activate()
$wsAssets.Range("A6").Value = "PO_LookupURL"
$wsAssets.Range("B6").Value = "PO_LookupURL"
$wsAssets.Range("A6").Select()

# ---------------------------------------------------------------------------
# Constants sheet - clear stale selection state
# ---------------------------------------------------------------------------
$wsConstants.Activate()
$wsConstants.Range("A1").Select()

# ---------------------------------------------------------------------------
# Leave the workbook with the Assets sheet active/selected, matching the
# final authoring session
# ---------------------------------------------------------------------------
$wsAssets.Activate()
$wsAssets.Range("A6").Select()
